# Applies "Some of the biggest recipe methods are joined" edit.
$wb = $excel.ActiveWorkbook

$wsMethods = $wb.Worksheets.Item("Methods Required")

# --- Update cell values on "Methods Required" sheet ---
# Row 23: mark D23 as done ("x")
$wsMethods.Range("D23").Value = "x"

# Row 27: mark C27 and D27 as done ("x")
$wsMethods.Range("C27").Value = "x"
$wsMethods.Range("D27").Value = "x"

# Row 28: mark C28 and D28 as done ("x")
$wsMethods.Range("C28").Value = "x"
$wsMethods.Range("D28").Value = "x"

# Row 29: mark C29 and D29 as done ("x")
# (C29 did not previously hold a value, so its format must be set to match
#  the rest of the column - centered horizontally and vertically, same as
#  the neighboring cells.)
$wsMethods.Range("C29").Value = "x"
$wsMethods.Range("C29").HorizontalAlignment = -4108
$wsMethods.Range("C29").VerticalAlignment = -4108
$wsMethods.Range("D29").Value = "x"

# Row 30: mark C30 and D30 as done ("x")
$wsMethods.Range("C30").Value = "x"
$wsMethods.Range("D30").Value = "x"

# Row 31: mark C31 and D31 as pending ("pend")
$wsMethods.Range("C31").Value = "pend"
$wsMethods.Range("D31").Value = "pend"

# --- Make "Methods Required" the active sheet (was "Hoja1"), with D31 as
#     the selected / active cell (was A10). ---
$wsMethods.Activate()
$wsMethods.Range("D31").Select() | Out-Null
